$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.508.37"
$ws.Range("E2").Value = "  +11.11%  "
$ws.Range("D3").Value = "3.488.85"
$ws.Range("E3").Value = "  +6.98%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "417.03"
$ws.Range("E5").Value = "  +4.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "122.93"
$ws.Range("E6").Value = "  +13.07%  "
$ws.Range("D7").Value = "3.482.29"
$ws.Range("E7").Value = "  +6.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.604"
$ws.Range("E8").Value = "  +4.51%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.686"
$ws.Range("E10").Value = "  +10.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.133"
$ws.Range("E11").Value = "  +39.45%  "
$ws.Range("E12").Value = "  +5.39%  "
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "4.035.71"
$ws.Range("E14").Value = "  +6.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.62"
$ws.Range("E15").Value = "  +4.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.01"
$ws.Range("E16").Value = "  +5.48%  "
$ws.Range("D17").Value = "3.483.09"
$ws.Range("E17").Value = "  +7.06%  "
$ws.Range("D18").Value = "63.226.07"
$ws.Range("E18").Value = "  +10.98%  "
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.04"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000145"
$ws.Range("E21").Value = "  +34.80%  "
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "83.07"
$ws.Range("E23").Value = "  +12.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "316.67"
$ws.Range("E24").Value = "  +7.52%  "
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.18"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "31.28"
$ws.Range("E27").Value = "  +11.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.74"
$ws.Range("E28").Value = "  +3.28%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +4.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.31"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("E32").Value = "  +2.77%  "
$ws.Range("E33").Value = "  +4.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.57"
$ws.Range("E34").Value = "  +20.34%  "
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "41.99"
$ws.Range("E36").Value = "  +4.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0490"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.30"
$ws.Range("E38").Value = "  +1.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.996"
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("E42").Value = "  +7.33%  "
$ws.Range("E43").Value = "  +4.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "135.62"
$ws.Range("E44").Value = "  -1.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.02"
$ws.Range("E45").Value = "  +1.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.283"
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.91"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.26"
$ws.Range("E48").Value = "  +2.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.06"
$ws.Range("E49").Value = "  -1.60%  "
$ws.Range("D50").Value = "3.826.55"
$ws.Range("E50").Value = "  +6.85%  "
$ws.Range("D51").Value = "2.187.76"
$ws.Range("E51").Value = "  +1.88%  "
